# Update the regression table in the active document.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Widen the second data column slightly (1283 -> 1344 twips = 64.15pt -> 67.2pt).
$t.Columns.Item(2).Width = 67.2

# Row 2: (Intercept) coefficients
$t.Cell(2,2).Range.Text = "0.207 ***"
$t.Cell(2,3).Range.Text = "0.207 **"
$t.Cell(2,4).Range.Text = "-0.813 ***"
$t.Cell(2,5).Range.Text = "-0.813 **"

# Row 3: (Intercept) standard errors
$t.Cell(3,2).Range.Text = "(0.057)   "
$t.Cell(3,3).Range.Text = "(0.077)  "
$t.Cell(3,4).Range.Text = "(0.173)   "
$t.Cell(3,5).Range.Text = "(0.275)  "

# Row 4: treatment coefficients
$t.Cell(4,2).Range.Text = "0.265 ** "
$t.Cell(4,3).Range.Text = "0.265 * "
$t.Cell(4,4).Range.Text = "0.727 ** "
$t.Cell(4,5).Range.Text = "0.727 * "

# Row 5: treatment standard errors
$t.Cell(5,2).Range.Text = "(0.083)   "
$t.Cell(5,3).Range.Text = "(0.103)  "
$t.Cell(5,4).Range.Text = "(0.236)   "
$t.Cell(5,5).Range.Text = "(0.287)  "

# Row 6: educated coefficients
$t.Cell(6,2).Range.Text = "0.128    "
$t.Cell(6,3).Range.Text = "0.128   "
$t.Cell(6,4).Range.Text = "0.366    "
$t.Cell(6,5).Range.Text = "0.366   "

# Row 7: educated standard errors
$t.Cell(7,2).Range.Text = "(0.088)   "
$t.Cell(7,3).Range.Text = "(0.072)  "
$t.Cell(7,4).Range.Text = "(0.250)   "
$t.Cell(7,5).Range.Text = "(0.227)  "

# Row 8: N
$t.Cell(8,2).Range.Text = "131        "
$t.Cell(8,3).Range.Text = "131       "
$t.Cell(8,4).Range.Text = "131       "

# Row 9: Clusters
$t.Cell(9,2).Range.Text = "        "
$t.Cell(9,3).Range.Text = "4       "
$t.Cell(9,5).Range.Text = "4       "

# Row 10: Adj. R2
$t.Cell(10,2).Range.Text = "0.074    "
$t.Cell(10,3).Range.Text = "       "
$t.Cell(10,4).Range.Text = "       "
